$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-03-t02-activity-added")

$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

$ws.Range("F1").Select()
